# Apply updated "results with all iterations" values to both worksheets.

$wb = $excel.ActiveWorkbook

$wsFull = $wb.Worksheets.Item("Full results")
$wsPlot = $wb.Worksheets.Item("For plotting")

# --- "Full results" sheet ---

# Row 2 (education / COMPLETE MODEL)
$wsFull.Range("H2").Value = 0.580913366522763
$wsFull.Range("I2").Value = 0.193973483700182
$wsFull.Range("O2").Value = 0.419292173058898

# Row 3 (education / CONDITIONAL MODEL)
$wsFull.Range("F3").Value = 0.567348953574101
$wsFull.Range("G3").Value = 0.239078265098996

# Row 4 (education / NULL MODEL)
$wsFull.Range("C4").Value = 0.603975554745592
$wsFull.Range("D4").Value = 0.396378393439593
$wsFull.Range("E4").Value = 1.00035394818519
$wsFull.Range("J4").Value = 0.396238145549144
$wsFull.Range("K4").Value = 0.238993672765317
$wsFull.Range("L4").Value = -0.0135596136053715
$wsFull.Range("M4").Value = 0.0230540275097533
$wsFull.Range("N4").Value = 0.225434059159945

# --- "For plotting" sheet ---

# Row 2 (Sibcorr)
$wsPlot.Range("C2").Value = 0.396238145549144
$wsPlot.Range("D2").Value = 0.337188975695219
$wsPlot.Range("E2").Value = 0.45528731540307

# Row 3 (IOLIB)
$wsPlot.Range("C3").Value = 0.225434059159945
$wsPlot.Range("D3").Value = 0.165366293271668
$wsPlot.Range("E3").Value = 0.285501825048223

# Row 4 (IORAD)
$wsPlot.Range("C4").Value = 0.419292173058898
$wsPlot.Range("D4").Value = 0.355470570506988
$wsPlot.Range("E4").Value = 0.483113775610807
